$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.172.98'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.271.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.63%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '298.20'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.64'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -5.59%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.494'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.48%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.13'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.57%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.33'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -7.66%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.65'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.624.12'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.280.86'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.775'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.170.70'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.36'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.98'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.61%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '232.66'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.61%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.78'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.02%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.23'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.66'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.07%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.74%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.68%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.63%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.74%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '16.03'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -8.95%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0988'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.29%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.59%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -7.28%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -7.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.955.64'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.60%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.89%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.44'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -6.56%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.78'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -4.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.495.88'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.46%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.90'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.85%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.75'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.92%  '
